# RegistracijaAgencije_Scenarij.xlsx - "Scenarij za registraciju agencije"
#
# Moves the scenario header fields (Naziv/Opis/Glavni tok/Preduvjeti/Posljedice)
# out of column C into column D as bold-label + plain-text rich strings, shifts
# the "Scenarij 2" title from A1 to B1, fixes the "Porvjera" -> "Provjera" typo,
# enables word-wrap on the long description/consequence cells, tweaks row
# heights / column widths, and updates the final selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 1: "Scenarij 2" moves from A1 to B1 (bold style carried over)
# ---------------------------------------------------------------------------
$titleValue = $ws.Range("A1").Value2
$ws.Range("B1").Value = $titleValue
$ws.Range("B1").Font.Bold = $true
$ws.Range("A1").ClearContents()

# ---------------------------------------------------------------------------
# Helper data: label (bold run) + rest-of-text (normal run) for each of the
# five header rows, moved from column C to column D.
# ---------------------------------------------------------------------------

# Row 3 : Naziv
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = "Naziv: Registracija agencije na stranicu"
$r = $ws.Range("D3").Characters(1, 5)
$r.Font.Bold = $true
$r.Font.Name = "Calibri"
$r.Font.Size = 11
$r2 = $ws.Range("D3").Characters(6, 35)
$r2.Font.Name = "Calibri"
$r2.Font.Size = 11
$ws.Rows(3).RowHeight = 15

# Row 4 : Opis (long description, wrapped)
$ws.Range("C4").ClearContents()
$ws.Range("D4").Value = "Opis: Ako agencija odluči da želi kreirati i reklamirati putovanja preko naše stranice mora da kreira račun. Vrši se unos podataka, te se provjerava njihova ispravnost. Ako je sve ispravno završeno je kreiranje računa. Kada se prijavi na račun pojavljuje se forma preko koje može da kreira putovanja tj.odabir destinacije,datum i vrijeme,hotel i let. Ako se odluči da kreira vrši se provjera kapaciteta i plaća određeni procenat. Agencija ne mora kreirati odma putovanja, može i da pogleda svoja prethodno kreirana.  "
$r = $ws.Range("D4").Characters(1, 4)
$r.Font.Bold = $true
$r.Font.Name = "Calibri"
$r.Font.Size = 11
$r2 = $ws.Range("D4").Characters(5, 513)
$r2.Font.Name = "Calibri"
$r2.Font.Size = 11
$ws.Range("C4").WrapText = $true
$ws.Range("D4").WrapText = $true
$ws.Rows(4).RowHeight = 157.2

# Row 5 : Glavni tok
$ws.Range("C5").ClearContents()
$ws.Range("D5").Value = "Glavni tok: Agencija kreira račun i kreira putovanje"
$r = $ws.Range("D5").Characters(1, 10)
$r.Font.Bold = $true
$r.Font.Name = "Calibri"
$r.Font.Size = 11
$r2 = $ws.Range("D5").Characters(11, 42)
$r2.Font.Name = "Calibri"
$r2.Font.Size = 11

# Row 6 : Preduvjeti
$ws.Range("C6").ClearContents()
$ws.Range("D6").Value = "Preduvjeti: Agencija mora imati račun na stranici"
$r = $ws.Range("D6").Characters(1, 10)
$r.Font.Bold = $true
$r.Font.Name = "Calibri"
$r.Font.Size = 11
$r2 = $ws.Range("D6").Characters(11, 39)
$r2.Font.Name = "Calibri"
$r2.Font.Size = 11

# Row 7 : Posljedice (wrapped)
$ws.Range("C7").ClearContents()
$ws.Range("D7").Value = "Posljedice: Ako odluči da otkaže putovanje procenti koje je već uplatila se ne vraćaju, jer su hotel i avio kompanija već izgubile klijente zbog rezervacije."
$r = $ws.Range("D7").Characters(1, 10)
$r.Font.Bold = $true
$r.Font.Name = "Calibri"
$r.Font.Size = 11
$r2 = $ws.Range("D7").Characters(11, 147)
$r2.Font.Name = "Calibri"
$r2.Font.Size = 11
$ws.Range("C7").WrapText = $true
$ws.Range("D7").WrapText = $true
$ws.Rows(7).RowHeight = 57.6

# ---------------------------------------------------------------------------
# Row 10: fix typo "Porvjera" -> "Provjera"
# ---------------------------------------------------------------------------
$ws.Range("D10").Value = "2. Provjera validnosti"

# ---------------------------------------------------------------------------
# Column widths (closest values reachable through the pixel-snapped
# ColumnWidth COM property)
# ---------------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 27.333333333333336
$ws.Columns("D").ColumnWidth = 40.33333333333333

# ---------------------------------------------------------------------------
# Final selection: D7, scrolled back to the top (no frozen topLeftCell offset)
# ---------------------------------------------------------------------------
$ws.Application.Goto($ws.Range("A1"))
$ws.Range("D7").Select()
